# Update COVID country stats, then re-sort by total cases (col B) descending,
# then update the "last updated" timestamp banner in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Apply the updated per-country figures (Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes).
#    Looked up by country name so the edit is correct regardless of each
#    country's current row position. Countries not listed here (e.g.
#    Barein, Uruguay, Albania, Taiwan) keep their own numbers unchanged -
#    they only end up on a different row because the re-sort below moves
#    Hungria/Afganistan past them once those two grow past their totals.

$updates = @(
    @("Estados Unidos", 367650, 646, 19810, 336897, 8983, 72, 10943),
    @("Australia",         5895,   0,  2432,   3417,   96,  1,    46),
    @("Pakistan",          3864,  98,   429,   3381,   28,  1,    54),
    @("Tailandia",         2258,  38,   824,   1407,   30,  1,    27),
    @("Lituania",           880,  37,     8,    857,   11,  0,    15),
    @("Hungria",            817,  73,    71,    699,   17,  9,    47),
    @("Bulgaria",           565,  16,    42,    501,   26,  0,    22),
    @("Afganistan",         423,  56,    18,    394,    0,  0,    11),
    @("Vietnam",            245,   0,   106,    139,    8,  0,     0)
)

foreach ($row in $updates) {
    $country = $row[0]
    # LookAt:=xlWhole (1) so e.g. "Estados Unidos" never matches the
    # "...Islas Virgenes de los Estados Unidos" row by substring.
    $found = $ws.Range("A4:A215").Find($country, [Type]::Missing, [Type]::Missing, 1)
    $r = $found.Row()
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# 2) Re-sort the data block (rows 4-215) by "Casos totales" (col B),
#    descending, keeping the relative order of ties (stable sort), same
#    as the dashboard does whenever the underlying numbers are refreshed.

$dataRange = $ws.Range("A4:H215")
$sortKey = $ws.Range("B4:B215")
$dataRange.Sort($sortKey, 2)

# 3) Update the "Datos actualizados..." banner.

$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 07:52"
